# Apply the "ADD NEW DAILY DATA ROW" update to the 相談件数 (consultation count) sheet.
# A new data row (date 2020-05-19 / serial 43970) is inserted as row 115, pushing the
# trailing footnote row down from 115 to 116. Print area / dimension / view selections
# are updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")

# Insert a new row above the old footer row (currently row 115) so the footer moves to 116,
# and the new data row takes its place at row 115, copying formatting from row 114 above it.
$ws.Rows.Item(115).Insert(-4121)  # xlShiftDown

# Row 115 needs to look like the preceding data rows (A: date, B/C counts, D/E totals).
$ws.Range("A115:E115").Value = $ws.Range("A114:E114").Value
$ws.Range("A115:E115").NumberFormat = $ws.Range("A114:E114").NumberFormat

# Fill in the new day's figures.
$ws.Range("A115").Value = 43970
$ws.Range("B115").Value = 206
$ws.Range("C115").Value = 38171
$ws.Range("D115").Value = 40
$ws.Range("E115").Value = 7682

# Update the print area / dimension to extend through the new last row (116).
$wb.Names.Item("_xlnm.Print_Area").RefersToR1C1 = "=相談件数!R1C1:R116C5"

# Update the frozen-pane / selection view state to match the shipped workbook.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C85").Select()
$ws.Application.ActiveWindow.FreezePanes = $false
$ws.Range("B2").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("D113").Select()
